$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.455362044514542
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 4.358119930609447

$ws.Range("B3").Value = 1.455362044514542
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 22.3905356188092
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 25.99591228164478

$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 22.3905356188092
$ws.Range("E4").Value = 10.19245300693656
$ws.Range("G4").Value = 37.52559925287081

$ws.Range("B5").Value = 3.286832544864788
$ws.Range("C5").Value = 1.655778082260271
$ws.Range("D5").Value = 3.537761648806719
$ws.Range("E5").Value = 0.4942365360607697
$ws.Range("G5").Value = 8.974608811992548

$ws.Range("B6").Value = 0.2917716402565462
$ws.Range("C6").Value = 0.04071648406533734
$ws.Range("D6").Value = 3.537761648806719
$ws.Range("E6").Value = 0.4942365360607697
$ws.Range("G6").Value = 4.364486309189372

$ws.Range("B7").Value = 3.286832544864788
$ws.Range("C7").Value = 1.655778082260271
$ws.Range("D7").Value = 0.1494219747398047
$ws.Range("E7").Value = 0.4942365360607697
$ws.Range("G7").Value = 5.586269137925634

$ws.Range("B8").Value = 3.286832544864788
$ws.Range("C8").Value = 1.655778082260271
$ws.Range("D8").Value = 0.7527432677738641
$ws.Range("E8").Value = 10.19245300693656
$ws.Range("G8").Value = 15.88780690183548

$ws.Range("B9").Value = 3.286832544864788
$ws.Range("C9").Value = 1.655778082260271
$ws.Range("D9").Value = 0.7527432677738641
$ws.Range("E9").Value = 0.4942365360607697
$ws.Range("G9").Value = 6.189590430959694

$ws.Range("B10").Value = 3.286832544864788
$ws.Range("C10").Value = 1.655778082260271
$ws.Range("D10").Value = 3.537761648806719
$ws.Range("E10").Value = 0.4942365360607697
$ws.Range("G10").Value = 8.974608811992548

$ws.Range("B11").Value = 3.286832544864788
$ws.Range("C11").Value = 1.655778082260271
$ws.Range("D11").Value = 0.7527432677738641
$ws.Range("E11").Value = 0.4942365360607697
$ws.Range("G11").Value = 6.189590430959694

$ws.Range("B12").Value = 1.455362044514542
$ws.Range("C12").Value = 0.306821227259698
$ws.Range("D12").Value = 0.1494219747398047
$ws.Range("E12").Value = 0.4942365360607697
$ws.Range("G12").Value = 2.405841782574814

$ws.Range("B13").Value = 3.286832544864788
$ws.Range("C13").Value = 1.655778082260271
$ws.Range("D13").Value = 22.3905356188092
$ws.Range("E13").Value = 0.4942365360607697
$ws.Range("G13").Value = 27.82738278199502

$ws.Range("B14").Value = 3.286832544864788
$ws.Range("C14").Value = 1.655778082260271
$ws.Range("D14").Value = 0.7527432677738641
$ws.Range("E14").Value = 0.4942365360607697
$ws.Range("G14").Value = 6.189590430959694

$ws.Range("B15").Value = 3.286832544864788
$ws.Range("C15").Value = 1.655778082260271
$ws.Range("D15").Value = 0.7527432677738641
$ws.Range("E15").Value = 0.4942365360607697
$ws.Range("G15").Value = 6.189590430959694

$ws.Range("B16").Value = 3.286832544864788
$ws.Range("C16").Value = 1.655778082260271
$ws.Range("D16").Value = 3.537761648806719
$ws.Range("E16").Value = 0.4942365360607697
$ws.Range("G16").Value = 8.974608811992548

$ws.Range("B17").Value = 1.455362044514542
$ws.Range("C17").Value = 1.655778082260271
$ws.Range("D17").Value = 0.1494219747398047
$ws.Range("E17").Value = 0.4942365360607697
$ws.Range("G17").Value = 3.754798637575387

$ws.Range("B18").Value = 0.1190320826869504
$ws.Range("C18").Value = 1.655778082260271
$ws.Range("D18").Value = 0.1494219747398047
$ws.Range("E18").Value = 0.4942365360607697
$ws.Range("G18").Value = 2.418468675747795

$ws.Range("B19").Value = 3.286832544864788
$ws.Range("C19").Value = 1.655778082260271
$ws.Range("D19").Value = 3.537761648806719
$ws.Range("E19").Value = 0.4942365360607697
$ws.Range("G19").Value = 8.974608811992548

$ws.Range("B20").Value = 3.286832544864788
$ws.Range("C20").Value = 1.655778082260271
$ws.Range("D20").Value = 0.7527432677738641
$ws.Range("E20").Value = 0.4942365360607697
$ws.Range("G20").Value = 6.189590430959694

$ws.Range("B21").Value = 3.286832544864788
$ws.Range("C21").Value = 1.655778082260271
$ws.Range("D21").Value = 0.1494219747398047
$ws.Range("E21").Value = 0.4942365360607697
$ws.Range("G21").Value = 5.586269137925634

$ws.Range("B22").Value = 3.286832544864788
$ws.Range("C22").Value = 1.655778082260271
$ws.Range("D22").Value = 3.537761648806719
$ws.Range("E22").Value = 0.4942365360607697
$ws.Range("G22").Value = 8.974608811992548

$ws.Range("B23").Value = 0.6606524410359556
$ws.Range("C23").Value = 1.655778082260271
$ws.Range("D23").Value = 3.537761648806719
$ws.Range("E23").Value = 0.4942365360607697
$ws.Range("G23").Value = 6.348428708163715

$ws.Range("B24").Value = 3.286832544864788
$ws.Range("C24").Value = 1.655778082260271
$ws.Range("D24").Value = 22.3905356188092
$ws.Range("E24").Value = 0.4942365360607697
$ws.Range("G24").Value = 27.82738278199502

$ws.Range("B25").Value = 3.286832544864788
$ws.Range("C25").Value = 1.655778082260271
$ws.Range("D25").Value = 0.7527432677738641
$ws.Range("E25").Value = 0.4942365360607697
$ws.Range("G25").Value = 6.189590430959694

$ws.Range("B26").Value = 3.286832544864788
$ws.Range("C26").Value = 1.655778082260271
$ws.Range("D26").Value = 22.3905356188092
$ws.Range("E26").Value = 0.4942365360607697
$ws.Range("G26").Value = 27.82738278199502
